$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45865, 19, 0.0, "27.07.202519"),
    @(45865, 20, 0.0, "27.07.202520"),
    @(45865, 21, 0.037, "27.07.202521"),
    @(45865, 22, 0.0, "27.07.202522"),
    @(45865, 23, 0.0, "27.07.202523"),
    @(45865, 24, 0.0, "27.07.202524"),
    @(45866, 1, 0.0, "28.07.20251"),
    @(45866, 2, 0.0, "28.07.20252"),
    @(45866, 3, 0.0, "28.07.20253"),
    @(45866, 4, 0.0, "28.07.20254"),
    @(45866, 5, 0.0, "28.07.20255"),
    @(45866, 6, 0.0, "28.07.20256"),
    @(45866, 7, 0.086, "28.07.20257"),
    @(45866, 8, 0.473, "28.07.20258"),
    @(45866, 9, 0.995, "28.07.20259"),
    @(45866, 10, 1.193, "28.07.202510"),
    @(45866, 11, 1.615, "28.07.202511"),
    @(45866, 12, 1.482, "28.07.202512"),
    @(45866, 13, 1.372, "28.07.202513"),
    @(45866, 14, 1.314, "28.07.202514"),
    @(45866, 15, 0.917, "28.07.202515"),
    @(45866, 16, 0.938, "28.07.202516"),
    @(45866, 17, 0.956, "28.07.202517"),
    @(45866, 18, 0.792, "28.07.202518"),
    @(45866, 19, 0.632, "28.07.202519"),
    @(45866, 20, 0.199, "28.07.202520"),
    @(45866, 21, 0.015, "28.07.202521"),
    @(45866, 22, 0.0, "28.07.202522"),
    @(45866, 23, 0.0, "28.07.202523"),
    @(45866, 24, 0.0, "28.07.202524"),
    @(45867, 1, 0.0, "29.07.20251"),
    @(45867, 2, 0.0, "29.07.20252"),
    @(45867, 3, 0.0, "29.07.20253"),
    @(45867, 4, 0.0, "29.07.20254"),
    @(45867, 5, 0.0, "29.07.20255"),
    @(45867, 6, 0.0, "29.07.20256"),
    @(45867, 7, 0.036, "29.07.20257"),
    @(45867, 8, 0.266, "29.07.20258"),
    @(45867, 9, 0.78, "29.07.20259"),
    @(45867, 10, 0.997, "29.07.202510"),
    @(45867, 11, 1.404, "29.07.202511"),
    @(45867, 12, 1.64, "29.07.202512"),
    @(45867, 13, 1.778, "29.07.202513"),
    @(45867, 14, 1.535, "29.07.202514"),
    @(45867, 15, 1.674, "29.07.202515"),
    @(45867, 16, 1.53, "29.07.202516"),
    @(45867, 17, 1.252, "29.07.202517"),
    @(45867, 18, 0.841, "29.07.202518"),
    @(45867, 19, 0.585, "29.07.202519"),
    @(45867, 20, 0.195, "29.07.202520"),
    @(45867, 21, 0.015, "29.07.202521"),
    @(45867, 22, 0.0, "29.07.202522"),
    @(45867, 23, 0.0, "29.07.202523"),
    @(45867, 24, 0.0, "29.07.202524"),
    @(45868, 1, 0.0, "30.07.20251"),
    @(45868, 2, 0.0, "30.07.20252"),
    @(45868, 3, 0.0, "30.07.20253"),
    @(45868, 4, 0.0, "30.07.20254"),
    @(45868, 5, 0.0, "30.07.20255"),
    @(45868, 6, 0.0, "30.07.20256"),
    @(45868, 7, 0.086, "30.07.20257"),
    @(45868, 8, 0.5, "30.07.20258"),
    @(45868, 9, 1.003, "30.07.20259"),
    @(45868, 10, 1.277, "30.07.202510"),
    @(45868, 11, 1.227, "30.07.202511"),
    @(45868, 12, 1.778, "30.07.202512"),
    @(45868, 13, 1.576, "30.07.202513"),
    @(45868, 14, 1.523, "30.07.202514"),
    @(45868, 15, 1.439, "30.07.202515"),
    @(45868, 16, 1.329, "30.07.202516"),
    @(45868, 17, 1.459, "30.07.202517"),
    @(45868, 18, 1.058, "30.07.202518"),
    @(45868, 19, 0.788, "30.07.202519"),
    @(45868, 20, 0.271, "30.07.202520"),
    @(45868, 21, 0.034, "30.07.202521"),
    @(45868, 22, 0.0, "30.07.202522"),
    @(45868, 23, 0.0, "30.07.202523"),
    @(45868, 24, 0.0, "30.07.202524"),
    @(45869, 1, 0.0, "31.07.20251"),
    @(45869, 2, 0.0, "31.07.20252"),
    @(45869, 3, 0.0, "31.07.20253"),
    @(45869, 4, 0.0, "31.07.20254"),
    @(45869, 5, 0.0, "31.07.20255"),
    @(45869, 6, 0.0, "31.07.20256"),
    @(45869, 7, 0.085, "31.07.20257"),
    @(45869, 8, 0.48, "31.07.20258"),
    @(45869, 9, 1.01, "31.07.20259"),
    @(45869, 10, 1.406, "31.07.202510"),
    @(45869, 11, 1.338, "31.07.202511"),
    @(45869, 12, 1.403, "31.07.202512"),
    @(45869, 13, 1.403, "31.07.202513"),
    @(45869, 14, 1.362, "31.07.202514"),
    @(45869, 15, 1.329, "31.07.202515"),
    @(45869, 16, 1.038, "31.07.202516"),
    @(45869, 17, 0.909, "31.07.202517"),
    @(45869, 18, 0.806, "31.07.202518"),
    @(45869, 19, 0.545, "31.07.202519"),
    @(45869, 20, 0.202, "31.07.202520"),
    @(45869, 21, 0.015, "31.07.202521"),
    @(45869, 22, 0.0, "31.07.202522"),
    @(45869, 23, 0.0, "31.07.202523"),
    @(45869, 24, 0.0, "31.07.202524"),
    @(45870, 1, 0.0, "01.08.20251"),
    @(45870, 2, 0.0, "01.08.20252"),
    @(45870, 3, 0.0, "01.08.20253"),
    @(45870, 4, 0.0, "01.08.20254"),
    @(45870, 5, 0.0, "01.08.20255"),
    @(45870, 6, 0.0, "01.08.20256"),
    @(45870, 7, 0.064, "01.08.20257"),
    @(45870, 8, 0.474, "01.08.20258"),
    @(45870, 9, 0.936, "01.08.20259"),
    @(45870, 10, 1.118, "01.08.202510"),
    @(45870, 11, 1.338, "01.08.202511"),
    @(45870, 12, 1.348, "01.08.202512"),
    @(45870, 13, 0.969, "01.08.202513"),
    @(45870, 14, 0.944, "01.08.202514"),
    @(45870, 15, 0.908, "01.08.202515"),
    @(45870, 16, 1.024, "01.08.202516"),
    @(45870, 17, 0.909, "01.08.202517"),
    @(45870, 18, 0.731, "01.08.202518"),
    @(45870, 19, 0.493, "01.08.202519"),
    @(45870, 20, 0.182, "01.08.202520"),
    @(45870, 21, 0.015, "01.08.202521"),
    @(45870, 22, 0.0, "01.08.202522"),
    @(45870, 23, 0.0, "01.08.202523"),
    @(45870, 24, 0.0, "01.08.202524"),
    @(45871, 1, 0.0, "02.08.20251"),
    @(45871, 2, 0.0, "02.08.20252"),
    @(45871, 3, 0.0, "02.08.20253"),
    @(45871, 4, 0.0, "02.08.20254"),
    @(45871, 5, 0.0, "02.08.20255"),
    @(45871, 6, 0.0, "02.08.20256"),
    @(45871, 7, 0.088, "02.08.20257"),
    @(45871, 8, 0.613, "02.08.20258"),
    @(45871, 9, 0.914, "02.08.20259"),
    @(45871, 10, 1.252, "02.08.202510"),
    @(45871, 11, 0.897, "02.08.202511"),
    @(45871, 12, 0.97, "02.08.202512"),
    @(45871, 13, 0.97, "02.08.202513"),
    @(45871, 14, 0.926, "02.08.202514"),
    @(45871, 15, 0.869, "02.08.202515"),
    @(45871, 16, 0.869, "02.08.202516"),
    @(45871, 17, 0.971, "02.08.202517"),
    @(45871, 18, 0.856, "02.08.202518"),
    @(45871, 19, 0.669, "02.08.202519"),
    @(45871, 20, 0.244, "02.08.202520"),
    @(45871, 21, 0.024, "02.08.202521"),
    @(45871, 22, 0.0, "02.08.202522"),
    @(45871, 23, 0.0, "02.08.202523"),
    @(45871, 24, 0.0, "02.08.202524"),
    @(45872, 1, 0.0, "03.08.20251"),
    @(45872, 2, 0.0, "03.08.20252"),
    @(45872, 3, 0.0, "03.08.20253"),
    @(45872, 4, 0.0, "03.08.20254"),
    @(45872, 5, 0.0, "03.08.20255"),
    @(45872, 6, 0.0, "03.08.20256"),
    @(45872, 7, 0.088, "03.08.20257"),
    @(45872, 8, 0.627, "03.08.20258"),
    @(45872, 9, 0.96, "03.08.20259"),
    @(45872, 10, 1.297, "03.08.202510"),
    @(45872, 11, 1.382, "03.08.202511"),
    @(45872, 12, 1.009, "03.08.202512"),
    @(45872, 13, 1.009, "03.08.202513"),
    @(45872, 14, 0.965, "03.08.202514"),
    @(45872, 15, 0.878, "03.08.202515"),
    @(45872, 16, 0.892, "03.08.202516"),
    @(45872, 17, 1.06, "03.08.202517"),
    @(45872, 18, 0.838, "03.08.202518"),
    @(45872, 19, 0.647, "03.08.202519")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}